$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold numeric-looking strings (e.g. "1.014", "29.139.45")
# that must stay text. A leading apostrophe forces Excel to treat the
# assigned value as text (quote-prefix) instead of coercing it to a number.

$ws.Range("D2").Value = "29.139.45"
$ws.Range("E2").Value = "  -1.54%  "
$ws.Range("D3").Value = "1.995.69"
$ws.Range("E3").Value = "  -0.51%  "
$ws.Range("D4").Value = "'1.014"
$ws.Range("E4").Value = "  +0.73%  "
$ws.Range("D5").Value = "'330.61"
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("D6").Value = "'1.012"
$ws.Range("E6").Value = "  +0.60%  "
$ws.Range("D7").Value = "'0.4981"
$ws.Range("E7").Value = "  -0.73%  "
$ws.Range("D8").Value = "'0.4197"
$ws.Range("E8").Value = "  -1.50%  "
$ws.Range("D9").Value = "'54.75"
$ws.Range("E9").Value = "  +1.72%  "
$ws.Range("D10").Value = "'0.08904"
$ws.Range("E10").Value = "  -3.05%  "
$ws.Range("D11").Value = "'1.096"
$ws.Range("E11").Value = "  -2.87%  "
$ws.Range("D12").Value = "'22.99"
$ws.Range("E12").Value = "  -2.36%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "'8.006"
$ws.Range("E13").Value = "  -1.55%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.975.88"
$ws.Range("E14").Value = "  +0.49%  "
$ws.Range("D15").Value = "'6.437"
$ws.Range("E15").Value = "  -1.77%  "
$ws.Range("D16").Value = "'1.014"
$ws.Range("E16").Value = "  +0.73%  "
$ws.Range("D17").Value = "'92.57"
$ws.Range("E17").Value = "  -3.67%  "
$ws.Range("D18").Value = "'0.00001107"
$ws.Range("E18").Value = "  -1.29%  "
$ws.Range("D19").Value = "'0.06768"
$ws.Range("E19").Value = "  +1.75%  "
$ws.Range("D20").Value = "'19.57"
$ws.Range("E20").Value = "  -1.52%  "
$ws.Range("E21").Value = "  +0.57%  "
$ws.Range("D22").Value = "'5.982"
$ws.Range("E22").Value = "  -0.24%  "
$ws.Range("D23").Value = "29.162.39"
$ws.Range("E23").Value = "  -1.45%  "
$ws.Range("D24").Value = "'11.99"
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("D25").Value = "'2.295"
$ws.Range("E25").Value = "  +1.25%  "
$ws.Range("D26").Value = "2.228.92"
$ws.Range("E26").Value = "  +0.52%  "
$ws.Range("D27").Value = "'20.87"
$ws.Range("E27").Value = "  +0.46%  "
$ws.Range("D28").Value = "'157.32"
$ws.Range("E28").Value = "  -1.03%  "
$ws.Range("D29").Value = "'6.323"
$ws.Range("E29").Value = "  -3.60%  "
$ws.Range("D30").Value = "'2.259"
$ws.Range("E30").Value = "  -3.75%  "
$ws.Range("D31").Value = "'127.40"
$ws.Range("E31").Value = "  -0.78%  "
$ws.Range("D32").Value = "'1.048"
$ws.Range("E32").Value = "  -1.13%  "
$ws.Range("D33").Value = "'0.09870"
$ws.Range("E33").Value = "  -0.87%  "
$ws.Range("D34").Value = "'1.530"
$ws.Range("E34").Value = "  -3.82%  "
$ws.Range("D35").Value = "'5.830"
$ws.Range("E35").Value = "  -0.73%  "
$ws.Range("D36").Value = "'3.749"
$ws.Range("E36").Value = "  -0.85%  "
$ws.Range("D37").Value = "'0.02422"
$ws.Range("E37").Value = "  -2.21%  "
$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D38").Value = "'9.164"
$ws.Range("E38").Value = "  -5.12%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "'1.316"
$ws.Range("E39").Value = "  +0.39%  "
$ws.Range("D40").Value = "'0.06407"
$ws.Range("E40").Value = "  +0.27%  "
$ws.Range("D41").Value = "'0.6502"
$ws.Range("E41").Value = "  -1.41%  "
$ws.Range("D42").Value = "'11.58"
$ws.Range("E42").Value = "  -1.97%  "
$ws.Range("D43").Value = "'0.1986"
$ws.Range("E43").Value = "  -4.54%  "
$ws.Range("E44").Value = "  +0.55%  "
$ws.Range("B45").Value = "WEMIXTOKEN"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").Value = "'1.366"
$ws.Range("E45").Value = "  +6.71%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.6219"
$ws.Range("E46").Value = "  -2.38%  "
$ws.Range("D47").Value = "'13.47"
$ws.Range("E47").Value = "  -0.47%  "
$ws.Range("D48").Value = "'2.186"
$ws.Range("E48").Value = "  -1.59%  "
$ws.Range("D49").Value = "'3.498"
$ws.Range("E49").Value = "  -1.02%  "
$ws.Range("D50").Value = "'0.00000000343"
$ws.Range("E50").Value = "  +6.51%  "
$ws.Range("D51").Value = "'2.199"
$ws.Range("E51").Value = "  +11.63%  "
